# Forms SOPs 1 & 2
# Rename the worksheet from the generic "Sheet1" to the SOP form code
# "F-SW-SD-05". Excel automatically keeps the sheet-scoped Print_Area
# defined name in sync with the new sheet name, but we set it explicitly
# too so the printed range stays pinned to A1:E13 regardless.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "F-SW-SD-05"

$ws.PageSetup.PrintArea = "`$A`$1:`$E`$13"
